$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the text of A20 (shared string "Test" -> "test") while keeping it referencing the same string
$ws.Range("A20").Value = "test"

# Fill in the KEYWORD column for this new pass case row
$ws.Range("D20").Value = "Keyword.properties"

# Update the active selection to A20, matching the saved view state
$ws.Range("A20").Select()
